$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "dateLastModified" -> "modified" (shared string used in A1)
$ws.Range("A1").Value = "modified"

# Update the active selection to A2 on this sheet
$ws.Range("A2").Select()
